$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CreateEmployee")
$ws.Range("D3").Value = "john doe"
